$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data

$ws.Range("D2").Value = "27.365.79"
$ws.Range("E2").Value = "  -3.84%  "

$ws.Range("D3").Value = "1.858.10"
$ws.Range("E3").Value = "  -4.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.94%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.43"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4540"
$ws.Range("E7").Value = "  -5.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3865"
$ws.Range("E8").Value = "  -5.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.45"
$ws.Range("E9").Value = "  -9.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07917"
$ws.Range("E10").Value = "  -7.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.018"
$ws.Range("E11").Value = "  -3.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.41"
$ws.Range("E12").Value = "  -4.65%  "

$ws.Range("D13").Value = "1.833.92"
$ws.Range("E13").Value = "  -8.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.916"
$ws.Range("E14").Value = "  -4.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.121"
$ws.Range("E15").Value = "  -6.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001032"
$ws.Range("E17").Value = "  -3.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "85.83"
$ws.Range("E18").Value = "  -5.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06508"
$ws.Range("E19").Value = "  -1.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.06"
$ws.Range("E20").Value = "  -7.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.533"
$ws.Range("E22").Value = "  -5.29%  "

$ws.Range("D23").Value = "27.365.96"
$ws.Range("E23").Value = "  -4.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.85"
$ws.Range("E24").Value = "  -5.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.276"
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").Value = "2.068.54"
$ws.Range("E26").Value = "  -7.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.61"
$ws.Range("E27").Value = "  -1.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("E28").Value = "  -2.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.067"
$ws.Range("E29").Value = "  -4.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.458"
$ws.Range("E30").Value = "  -6.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.73"
$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.486"
$ws.Range("E32").Value = "  +2.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09299"
$ws.Range("E33").Value = "  -3.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9333"
$ws.Range("E34").Value = "  -5.41%  "

$ws.Range("E35").Value = "  -2.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.262"
$ws.Range("E36").Value = "  -6.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02237"
$ws.Range("E37").Value = "  -4.30%  "

$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05995"
$ws.Range("E39").Value = "  -3.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.242"
$ws.Range("E40").Value = "  -10.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  -0.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5910"
$ws.Range("E42").Value = "  -5.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1885"
$ws.Range("E43").Value = "  -1.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.11"
$ws.Range("E44").Value = "  -9.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.263"
$ws.Range("E45").Value = "  -5.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5622"
$ws.Range("E46").Value = "  -5.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.93"
$ws.Range("E47").Value = "  -7.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.373"
$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.922"
$ws.Range("E49").Value = "  -6.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06777"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.31"
$ws.Range("E51").Value = "  -2.34%  "

